$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Registration")

# --- Update Sheet1 (Registration): drop firstName/lastName/zipCode columns ---
# NOTE: Range.Hyperlinks.Delete() removes every hyperlink on the sheet, so
# clear them all up-front and re-add the ones that should survive.
$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Range("A3:D3").ClearContents()
$ws1.Range("B1:D2").Clear()

$ws1.Range("A1").Value = "emailAddress"
$ws1.Range("A2").Value = "abctest@test.com"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:abctest@test.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "abctest@test.com")
$ws1.Range("A2").Style = "Hyperlink"

# --- Add Sheet2 (Login) right after Registration ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Login"

$ws2.Range("A2").Value = "bbb@bbb.com"
$ws2.Range("B2").Value = "bbbbb"
$ws2.Range("A1").Value = "userName"
$ws2.Range("B1").Value = "passWord"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:bbb@bbb.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "bbb@bbb.com")
$ws2.Range("A2").Style = "Hyperlink"

$ws2.Columns.Item(1).AutoFit()
$ws2.Range("E2").Select()
$ws2.Activate()

# --- Finish Sheet1 (Registration): set new last-row email + hyperlink ---
$ws1.Range("A3").Value = "baa@bbb.com"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:baa@bbb.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "baa@bbb.com")
$ws1.Range("A3").Style = "Hyperlink"

$ws1.Range("A11").Select()
$ws1.Activate()
